# Upload new version with timestamp
# - Inserts two new shortage items (alphabetically placed) into the report:
#     "HYDROFERRIN 50MG/ML ORAL DROPS 30 ML" and "KERELLA LOTION 30 ML"
#   between "HUSH SACHET" (row 28) and "LICID LOTION 30 ML" (old row 29).
# - Renumbers the following rows' item index (column A) by +2.
# - Updates the grand total (old P49) to include the two new prices.
# - Updates the footer timestamp.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert two blank rows right above the current "LICID LOTION 30 ML" row (row 29).
#    Everything currently on rows 29-50 shifts down to rows 31-52.
$ws.Rows("29:30").Insert()

# 2) Copy the formatting (styles, number formats, merges) of the row that is now
#    row 31 (the original row 29 template) down into the two freshly inserted rows.
#    NOTE: restrict to columns A:Q (the table's actual extent) - copying/pasting
#    whole rows would stamp formatting all the way out to column XFD and bloat
#    the sheet's used range.
$ws.Range("A31:Q31").Copy()
$ws.Range("A29:Q29").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A30:Q30").PasteSpecial(-4122)   # xlPasteFormats

# Re-establish the merged cell ranges for the two new rows (same layout as every
# other item row: A:B, C:G, H:K, L:M, N:O merged; P and Q stand alone).
$ws.Range("A29:B29").Merge()
$ws.Range("C29:G29").Merge()
$ws.Range("H29:K29").Merge()
$ws.Range("L29:M29").Merge()
$ws.Range("N29:O29").Merge()

$ws.Range("A30:B30").Merge()
$ws.Range("C30:G30").Merge()
$ws.Range("H30:K30").Merge()
$ws.Range("L30:M30").Merge()
$ws.Range("N30:O30").Merge()

# Match the row heights used throughout the rest of the table.
$ws.Rows("29:29").RowHeight = 25.5
$ws.Rows("30:30").RowHeight = 24.75

# 3) Fill in the values for the two new item rows.
$ws.Range("A29").Value = 23
$ws.Range("C29").Value = "HYDROFERRIN 50MG/ML ORAL DROPS 30 ML"
$ws.Range("H29").Value = "0:0"
$ws.Range("L29").Value = 1
$ws.Range("N29").Value = "44.00"
$ws.Range("P29").Value = "44.0000"
$ws.Range("Q29").Value = "1:0"

$ws.Range("A30").Value = 24
$ws.Range("C30").Value = "KERELLA LOTION 30 ML"
$ws.Range("H30").Value = "7:0"
$ws.Range("L30").Value = 1
$ws.Range("N30").Value = "31.00"
$ws.Range("P30").Value = "31.0000"
$ws.Range("Q30").Value = "1:0"

# 4) Renumber the item index column for all the original rows that shifted down
#    (now rows 31-50), so the sequence stays 1..44 with no gaps.
for ($r = 31; $r -le 50; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 6
}

# 5) Update the grand total (now on row 51) to account for the two new prices.
$ws.Range("P51").Value = 2800.5

# 6) Update the footer timestamp (now on row 52).
$ws.Range("A52").Value = "Sunday, 5 October, 2025 1:11 PM"
